$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# The test case with idCaso "3" (row 4, a blank placeholder row) is removed
# entirely; every row below it shifts up by one, which also drops the now-
# unused "3" entry from the shared string table.
$ws.Rows.Item(4).Delete()

# The case that used to hold "OSVPPRU01" (row 6 before the shift, row 5
# after it) gets a new "usuario" value.
$ws.Range("G5").Value = "chipote87"

# Reflect the author's final view/selection state in the saved file.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 6
$ws.Range("H10").Select()
